# Apply updated trial data (columns F:S, rows 2:29) for sheet "Sheet1".
# Columns A:E (subject_id..trial_total index, block info) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(189, 'bedrooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_0eflx.png', 76.05128205128206, 53.53846153846154, 64.7948717948718, 39, 6, 6, 6),
    @(190, 'bedrooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_5yhyk.png', 46.375, 31.325, 38.85, 40, 2, 2, 2),
    @(191, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_qgbyn.png', 65.08108108108108, 40.10810810810811, 52.5945945945946, 37, 4, 4, 4),
    @(192, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_h0hbk.png', 86.80952380952381, 69.19047619047619, 78, 42, 9, 9, 9),
    @(193, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_5m6x4.png', 80.23076923076923, 58.41025641025641, 69.32051282051282, 39, 7, 7, 7),
    @(194, 'bedrooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_zgg62.png', 82.1842105263158, 63.52631578947368, 72.85526315789474, 38, 8, 8, 8),
    @(195, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_a9acb.png', 77.11428571428571, 58.42857142857143, 67.77142857142857, 35, 7, 7, 7),
    @(196, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_x0u5z.png', 92, 78.16216216216216, 85.08108108108108, 37, 10, 10, 10),
    @(197, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_wyctg.png', 33.44736842105263, 11.39473684210526, 22.42105263157895, 38, 1, 1, 1),
    @(198, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_fqgem.png', 80.75, 61.475, 71.1125, 40, 8, 8, 8),
    @(199, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_v8dra.png', 61.77272727272727, 38.79545454545455, 50.28409090909091, 44, 3, 3, 3),
    @(200, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_uxxo0.png', 71.74418604651163, 48.44186046511628, 60.09302325581395, 43, 5, 5, 5),
    @(201, 'bedrooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_zv0dq.png', 76.86842105263158, 52.71052631578947, 64.78947368421052, 38, 6, 6, 6),
    @(202, 'bedrooms', $null, $null, 'catch', 'f', 'stimuli/catch_09.jpg', $null, $null, $null, $null, $null, $null, $null),
    @(203, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_ybbmx.png', 55.24324324324324, 36.75675675675676, 46, 37, 3, 3, 3),
    @(204, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_t2ioc.png', 88.1891891891892, 74.05405405405405, 81.12162162162161, 37, 10, 10, 10),
    @(205, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_th7xh.png', 82.35897435897436, 65.53846153846153, 73.94871794871796, 39, 8, 8, 8),
    @(206, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_71mhq.png', 69.34210526315789, 47.02631578947368, 58.18421052631579, 38, 5, 5, 5),
    @(207, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_oou46.png', 75.70270270270271, 54.86486486486486, 65.28378378378379, 37, 6, 6, 6),
    @(208, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_2js6m.png', 40.02777777777778, 20.88888888888889, 30.45833333333334, 36, 2, 2, 2),
    @(209, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_okvvw.png', 50.58333333333334, 32.11111111111111, 41.34722222222223, 36, 2, 2, 2),
    @(210, 'bedrooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_bklr1.png', 86.54761904761905, 67.73809523809524, 77.14285714285714, 42, 9, 9, 9),
    @(211, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_2pk6v.png', 85.08108108108108, 66.16216216216216, 75.62162162162161, 37, 9, 9, 9),
    @(212, 'bedrooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_i7vab.png', 86.4, 67.8, 77.1, 35, 9, 9, 9),
    @(213, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_rvssl.png', 74.25, 54.33333333333334, 64.29166666666667, 36, 6, 6, 6),
    @(214, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_le8uf.png', 12.88888888888889, 9.222222222222221, 11.05555555555556, 36, 1, 1, 1),
    @(215, 'bedrooms', 'bedrooms', $null, 'new', 'f', 'stimuli/img_6ddrx.png', 82.2, 63.68571428571428, 72.94285714285715, 35, 8, 8, 8),
    # Note: p_perceptual (S) intentionally stays 2 while p_typicality/p_conceptual become 4, per source data.
    @(216, 'bedrooms', 'bedrooms', 'target', 'old', 'j', 'stimuli/img_bj2gr.png', 65.25, 44.8, 55.025, 40, 4, 4, 2)
)

$nRows = 28
$nCols = 14
$data = New-Object 'object[,]' $nRows,$nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $data[$i,$j] = $rows[$i][$j]
    }
}

$ws.Range("F2:S29").Value2 = $data
